$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.974.67"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.754.11"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'577.49"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "'157.97"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "'0.111"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "'5.80"
$ws.Range("E10").Value = "  -13.88%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "3.239.51"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'26.76"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "63.868.73"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "'0.0000153"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "2.755.64"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'359.08"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "'0.553"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'66.12"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'8.46"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  +3.92%  "
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").Value = "'168.81"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "'20.29"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'4.16"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.13"
$ws.Range("E40").Value = "  +11.07%  "
$ws.Range("D41").Value = "'330.23"
$ws.Range("E41").Value = "  -4.28%  "
$ws.Range("D42").Value = "'39.35"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "'21.58"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "'21.65"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0256"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.633"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "'135.94"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'11.04"
$ws.Range("E51").Value = "  +0.65%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D51").ClearFormats()
